$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95; this shifts existing rows 95-171 down to 96-172.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new record.
$ws.Cells.Item(95, 1).Value = 3
$ws.Cells.Item(95, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44651
$ws.Cells.Item(95, 5).Value = 5
$ws.Cells.Item(95, 6).Value = 100112030
$ws.Cells.Item(95, 7).Value = "Poroto granado"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 73
$ws.Cells.Item(95, 11).Value = 21000
$ws.Cells.Item(95, 12).Value = 22000
$ws.Cells.Item(95, 13).Value = 21479
$ws.Cells.Item(95, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(95, 15).Value = "Provincia de Talca"
$ws.Cells.Item(95, 16).Value = 859
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
